$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.268.27'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '3.294.97'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.46'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.293.16'
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("D9").ClearFormats()
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.48'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000243'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.36'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '3.838.11'
$ws.Range("E15").Value = '  +1.33%  '
$ws.Range("E16").Value = '  +0.83%  '
$ws.Range("D17").Value = '3.301.92'
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("D18").Value = '63.358.47'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.05'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.89'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.90'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.78'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.96'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.75'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.10'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  +2.99%  '
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.09'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.17'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = '0.0₃0727'
$ws.Range("E38").Value = '  +2.69%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0397'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.29%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '3.111.50'
$ws.Range("E40").Value = '  +4.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '425.65'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.119'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +8.11%  '
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("E44").Value = '  -2.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.260'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.53%  '
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.07'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +8.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.08'
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.07'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.28'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.06%  '
